# timesheet.xlsx edit:
#   Insert two new timesheet rows (2026-02-12, 2026-02-15) right before the
#   "Total Duration" summary row, which pushes that summary row from row 24
#   down to row 26, and bump the grand total from "29 Hours" to "30 Hours".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the old row 24 ("Total Duration:" / "29
# Hours"). Excel-style row insert shifts that row down to row 26 and carries
# its existing values/formatting with it; the freshly inserted rows 24-25
# pick up the formatting of the row below, which already matches the style
# used by every other data row in the sheet.
$ws.Rows("24:25").Insert()

# ---- Row 24: 2026-02-12, 20:45:46 -> 21:04:21, 0.31 Hours ----
$ws.Range("B24").Value = "20:45:46"
$ws.Range("C24").Value = "21:04:21"
$ws.Range("D24").Value = "0.31 Hours"

# ---- Row 25: 2026-02-15, 17:02:42 -> 17:57:29, 0.91 Hours ----
$ws.Range("B25").Value = "17:02:42"
$ws.Range("C25").Value = "17:57:29"
$ws.Range("D25").Value = "0.91 Hours"

# The Date column values ("2026-02-12", "2026-02-15") look like dates, and a
# plain .Value assignment would auto-convert them into date serial numbers
# instead of keeping them as literal text (every other Date cell in this
# sheet stores the date as plain text, not a real Excel date). Stage each
# value through a scratch cell that is explicitly formatted as Text first,
# so the literal string is preserved, then copy it into place with a
# values-only paste.
$scratch = $ws.Range("Z1")

$scratch.NumberFormat = "@"
$scratch.Value = "2026-02-12"
$scratch.Copy()
$ws.Range("A24").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "2026-02-15"
$scratch.Copy()
$ws.Range("A25").PasteSpecial(-4163)

$scratch.Clear()

# Re-apply the normal (non-Text) cell formatting used throughout the sheet to
# the two date cells, since the values-only paste above left them without
# the shared "Times New Roman" style that the rest of the column uses.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 26 (formerly row 24): update the grand total ----
$ws.Range("D26").Value = "30 Hours"
